# The presentation ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> was "Office Theme" (used by the Notes Master)
#   ppt/theme/theme2.xml  -> was "Integral"     (used by the Slide Master / all slides)
#
# The authored change swaps the two themes' color schemes so that the theme
# driving the slides/Slide Master (theme2.xml) becomes the stock "Office"
# palette, while the Notes Master's theme effectively carries the former
# "Integral" palette.
#
# PowerPoint's object model exposes the *active* (Slide Master) theme's
# 12 DrawingML theme colors through Slide.ThemeColorScheme (Dark1, Light1,
# Dark2, Light2, Accent1-6, Hyperlink, FollowedHyperlink, in that order) -
# every slide shares the same Slide Master / theme part, so editing it via
# any one slide updates ppt/theme/theme2.xml for the whole deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette = the stock Office theme colors (RGB() = r + g*256 + b*65536)
$tcs.Item(1).RGB  = 0          # Dark1        -> 000000
$tcs.Item(2).RGB  = 16777215   # Light1       -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # Dark2        -> 44546A
$tcs.Item(4).RGB  = 15132391   # Light2       -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # Accent1      -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # Accent2      -> ED7D31
$tcs.Item(7).RGB  = 10855845   # Accent3      -> A5A5A5
$tcs.Item(8).RGB  = 49407      # Accent4      -> FFC000
$tcs.Item(9).RGB  = 12874308   # Accent5      -> 4472C4
$tcs.Item(10).RGB = 4697456    # Accent6      -> 70AD47
$tcs.Item(11).RGB = 12673797   # Hyperlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # FollowedHyperlink -> 954F72
